$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.529.66'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '1.826.84'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5161'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3891'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08433'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.00%  '
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.542'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '1.827.08'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.87%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001135'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06624'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.088'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '28.570.88'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  +2.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.273'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').Value = '2.037.60'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.408'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.83'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1097'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.099'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.07760'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.45%  '
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('B36').Value = 'Algorand'
$ws.Range('C36').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2229'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02382'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.265'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.751'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6370'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('E41').Value = '  +1.87%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6073'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.780'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.46'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.40%  '
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.207'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06984'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.68'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.04%  '
